# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# 1. Rename "Requested quantity" header on the Weekly Quantity sheet to Weekly_PO_Qty
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename "Requested quantity" header on the Monthly Trend sheet to Monthly_PO_Qty
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add a new "PO Forecast" sheet after "Monthly Trend" with forecast data
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "PO Forecast"

$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"
$ws.Range("A2").Value = 45207.99999999999
$ws.Range("B2").Value = 34
$ws.Range("C2").Value = -44.02191452654233
$ws.Range("D2").Value = 109.8533085464296
$ws.Range("A3").Value = 45214.99999999999
$ws.Range("B3").Value = 36
$ws.Range("C3").Value = -43.49207903812316
$ws.Range("D3").Value = 113.1987009475161
$ws.Range("A4").Value = 45221.99999999999
$ws.Range("B4").Value = 37
$ws.Range("C4").Value = -41.62323773759552
$ws.Range("D4").Value = 113.0924053825289
$ws.Range("A5").Value = 45228.99999999999
$ws.Range("B5").Value = 39
$ws.Range("C5").Value = -38.32138073476487
$ws.Range("D5").Value = 122.5684055146223
$ws.Range("A6").Value = 45235.99999999999
$ws.Range("B6").Value = 41
$ws.Range("C6").Value = -46.01196895418298
$ws.Range("D6").Value = 120.7098980346166
$ws.Range("A7").Value = 45242.99999999999
$ws.Range("B7").Value = 42
$ws.Range("C7").Value = -35.16440707324321
$ws.Range("D7").Value = 119.3763589988502
$ws.Range("A8").Value = 45249.99999999999
$ws.Range("B8").Value = 44
$ws.Range("C8").Value = -32.87488955086462
$ws.Range("D8").Value = 119.2290992835188
$ws.Range("A9").Value = 45256.99999999999
$ws.Range("B9").Value = 46
$ws.Range("C9").Value = -27.75405349962437
$ws.Range("D9").Value = 123.4295696803812
$ws.Range("A10").Value = 45263.99999999999
$ws.Range("B10").Value = 47
$ws.Range("C10").Value = -29.51464388596382
$ws.Range("D10").Value = 124.3332381019845
$ws.Range("A11").Value = 45270.99999999999
$ws.Range("B11").Value = 49
$ws.Range("C11").Value = -31.59432355856576
$ws.Range("D11").Value = 126.5645834503015
$ws.Range("A12").Value = 45277.99999999999
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = -27.7757045965584
$ws.Range("D12").Value = 133.3053607462177
$ws.Range("A13").Value = 45298.99999999999
$ws.Range("B13").Value = 55
$ws.Range("C13").Value = -23.53463627187392
$ws.Range("D13").Value = 128.7169956484204
$ws.Range("A14").Value = 45305.99999999999
$ws.Range("B14").Value = 57
$ws.Range("C14").Value = -23.83698998175128
$ws.Range("D14").Value = 129.2288099128286
$ws.Range("A15").Value = 45326.99999999999
$ws.Range("B15").Value = 62
$ws.Range("C15").Value = -18.59482204134401
$ws.Range("D15").Value = 142.8278676041376
$ws.Range("A16").Value = 45333.99999999999
$ws.Range("B16").Value = 63
$ws.Range("C16").Value = -11.92810926756419
$ws.Range("D16").Value = 144.707848878853
$ws.Range("A17").Value = 45340.99999999999
$ws.Range("B17").Value = 65
$ws.Range("C17").Value = -14.76787079828115
$ws.Range("D17").Value = 143.8404263219089
$ws.Range("A18").Value = 45347.99999999999
$ws.Range("B18").Value = 67
$ws.Range("C18").Value = -14.71620855644351
$ws.Range("D18").Value = 143.4564553337471
$ws.Range("A19").Value = 45354.99999999999
$ws.Range("B19").Value = 68
$ws.Range("C19").Value = -18.65709000801304
$ws.Range("D19").Value = 142.7439876222278
$ws.Range("A20").Value = 45361.99999999999
$ws.Range("B20").Value = 70
$ws.Range("C20").Value = -2.606350181348718
$ws.Range("D20").Value = 147.7256031565358
$ws.Range("A21").Value = 45368.99999999999
$ws.Range("B21").Value = 71
$ws.Range("C21").Value = -8.050195759012286
$ws.Range("D21").Value = 151.0978615357063
$ws.Range("A22").Value = 45375.99999999999
$ws.Range("B22").Value = 73
$ws.Range("C22").Value = -5.863650860818029
$ws.Range("D22").Value = 149.6744623994887
$ws.Range("A23").Value = 45382.99999999999
$ws.Range("B23").Value = 75
$ws.Range("C23").Value = -0.4204351707952637
$ws.Range("D23").Value = 154.4564946209299
$ws.Range("A24").Value = 45389.99999999999
$ws.Range("B24").Value = 76
$ws.Range("C24").Value = -2.564951989611826
$ws.Range("D24").Value = 152.9999541601076
$ws.Range("A25").Value = 45396.99999999999
$ws.Range("B25").Value = 78
$ws.Range("C25").Value = 2.22223922739349
$ws.Range("D25").Value = 156.5000817674527
$ws.Range("A26").Value = 45410.99999999999
$ws.Range("B26").Value = 81
$ws.Range("C26").Value = 7.105720819869746
$ws.Range("D26").Value = 160.5863324549938
$ws.Range("A27").Value = 45424.99999999999
$ws.Range("B27").Value = 84
$ws.Range("C27").Value = 8.374430294725972
$ws.Range("D27").Value = 161.613086968642
$ws.Range("A28").Value = 45431.99999999999
$ws.Range("B28").Value = 86
$ws.Range("C28").Value = 12.25110697114849
$ws.Range("D28").Value = 165.2912189523842
$ws.Range("A29").Value = 45438.99999999999
$ws.Range("B29").Value = 88
$ws.Range("C29").Value = 11.43743943118676
$ws.Range("D29").Value = 171.3646149900763
$ws.Range("A30").Value = 45459.99999999999
$ws.Range("B30").Value = 93
$ws.Range("C30").Value = 9.139413429387522
$ws.Range("D30").Value = 166.7748794274612
$ws.Range("A31").Value = 45466.99999999999
$ws.Range("B31").Value = 94
$ws.Range("C31").Value = 16.8403475342181
$ws.Range("D31").Value = 170.0345370730795
$ws.Range("A32").Value = 45473.99999999999
$ws.Range("B32").Value = 96
$ws.Range("C32").Value = 17.06738905907024
$ws.Range("D32").Value = 173.3225149533166
$ws.Range("A33").Value = 45480.99999999999
$ws.Range("B33").Value = 97
$ws.Range("C33").Value = 21.29709866357589
$ws.Range("D33").Value = 174.5119777502521
$ws.Range("A34").Value = 45487.99999999999
$ws.Range("B34").Value = 99
$ws.Range("C34").Value = 22.74079084719952
$ws.Range("D34").Value = 180.1937651142229
$ws.Range("A35").Value = 45494.99999999999
$ws.Range("B35").Value = 101
$ws.Range("C35").Value = 14.44256717310279
$ws.Range("D35").Value = 173.7469780789223
$ws.Range("A36").Value = 45501.99999999999
$ws.Range("B36").Value = 102
$ws.Range("C36").Value = 24.54801457444553
$ws.Range("D36").Value = 177.6282163120631
$ws.Range("A37").Value = 45508.99999999999
$ws.Range("B37").Value = 104
$ws.Range("C37").Value = 29.51701251330278
$ws.Range("D37").Value = 184.7256240868652
$ws.Range("A38").Value = 45515.99999999999
$ws.Range("B38").Value = 106
$ws.Range("C38").Value = 28.96309793921921
$ws.Range("D38").Value = 184.4119623341998
$ws.Range("A39").Value = 45522.99999999999
$ws.Range("B39").Value = 107
$ws.Range("C39").Value = 24.15445465532171
$ws.Range("D39").Value = 180.0982064238989

# Match header styling (bold, centered, bordered) used on the other sheets
$wsWeekly.Range("B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Match the date-formatted style used for the date column on the other sheets
$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A39").PasteSpecial(-4122)

$excel.CutCopyMode = $false

Write-Output "PO Forecast sheet added and headers updated"
